# Updated the comment about Ye helping Tabitha with OpenDS
#
# Original bullet text (User Support / TRANSCEND):
#   "TRANSCEND: Due to access-related issues, Abe could not walk Tabitha
#    through setting up users in OpenDS. Ye will take over now that he's
#    back."
#
# New bullet text:
#   "TRANSCEND: Ye walked Tabitha through creating her accounts."

$d = $word.ActiveDocument
$apos = [char]0x2019

# 1. Drop the leading "Due to access-related issues, " clause entirely -
#    it lived in its own run (Arial, cs="Arial") right after "TRANSCEND: ".
$d.Content.Find.Execute( `
    "Due to access-related issues, ", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 2) | Out-Null

# 2. Replace the remainder of the sentence. Doing this as a single
#    Find/Replace on the whole span would collapse everything into one run,
#    so instead: first overwrite the "Abe ... back." span with "Ye " using
#    Range.Text (which keeps the original run's Arial/no-cs formatting),
#    then grow the sentence with InsertAfter calls, nudging Font.Name back
#    to "Arial" after each insert so every new chunk gets its own run with
#    the matching rFonts instead of inheriting bare/default formatting.
$oldTail = "Abe could not walk Tabitha through setting up users in OpenDS." + `
    " Ye will take over now that he" + $apos + "s back."

$rng = $d.Content
$rng.Find.Execute($oldTail, $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null

$rng.Text = "Ye "

$rng.Collapse(0)
$rng.InsertAfter("walked Tabitha through creating her accounts")
$rng.Font.Name = "Arial"

$rng.Collapse(0)
$rng.InsertAfter(".")
$rng.Font.Name = "Arial"
